# This workbook's weekly price rows (2..46) get their "observation" columns
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Origen, Precio $/Kg) reshuffled across rows -- i.e. each row ends up showing
# the values that used to belong to a different row, while the descriptive
# columns (Mercado, Region, Categoria, Variedad, Calidad, Unidad, Kg o
# Unidades, Clasificacion) stay put. Row 45 keeps its own values (no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> source row (whose D/J/K/L/M/O/P values should land on "row")
$map = @{
    2  = 12
    3  = 15
    4  = 26
    5  = 8
    6  = 19
    7  = 32
    8  = 3
    9  = 37
    10 = 5
    11 = 25
    12 = 22
    13 = 36
    14 = 35
    15 = 10
    16 = 44
    17 = 24
    18 = 38
    19 = 13
    20 = 9
    21 = 28
    22 = 34
    23 = 29
    24 = 46
    25 = 6
    26 = 41
    27 = 2
    28 = 30
    29 = 7
    30 = 17
    31 = 16
    32 = 14
    33 = 20
    34 = 39
    35 = 27
    36 = 43
    37 = 4
    38 = 33
    39 = 40
    40 = 31
    41 = 11
    42 = 23
    43 = 18
    44 = 42
    45 = 45
    46 = 21
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for the columns being shuffled before any
# writes happen, so that later writes don't clobber values still needed as a
# source for other rows.
$snapshot = @{}
for ($r = 2; $r -le 46; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($r in $map.Keys) {
    $src = $map[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
